# Automatic update of files.
#
# The "Förändrad" column (C) holds a date serial that the logging sheet
# bumps by one day (2025-02-03 -> 2025-02-04, serials 45691 -> 45692) for
# every data row (rows 2-36) on the "Avverkningsanmälningar" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

for ($row = 2; $row -le 36; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    $cell.Value = 45692
}
